# Apply the edits described in the commit "Add files via upload":
#  - round two of the computed percentage-style values in row 2 and row 3
#    (columns F/G) down to whole numbers and switch B2:G3 to an integer
#    number format
#  - clear out the "test" / 45,56,4,55,99,31 sample row (row 4) so the
#    "test" shared string is no longer referenced
#  - swap the cell formatting that had been accidentally applied to A31 and
#    A32
#  - remove a spare blank row (the sheet shrinks from 109 rows to 108 rows)
#  - move the active selection to G3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 / Row 3: replace the long decimal values with their rounded
#     whole-number equivalents, then format B2:G3 as integers -----------
$ws.Range("F2").Value = 72
$ws.Range("G2").Value = 46
$ws.Range("F3").Value = 57
$ws.Range("G3").Value = 35
$ws.Range("B2:G3").NumberFormat = "0"

# --- Row 4: clear the leftover "test" sample row ------------------------
$ws.Range("A4:G4").ClearContents()

# --- Swap the formatting of A31 and A32 ---------------------------------
$ws.Range("A32").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("A32").PasteSpecial(-4122)

# --- Remove a spare blank row (row 85), shifting rows below it up -------
$ws.Rows(85).Delete()

# --- Move / record the active selection on the sheet --------------------
$ws.Range("G3").Select()
